$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for rows 3-20 (columns A-F), per the updated stats data
$data = @(
    @(901,  16, 15, 45, 60, 60),   # row 3
    @(401,  9,  48, 67, 75, 45),   # row 4
    @(201,  9,  30, 15, 45, 30),   # row 5
    @(701,  3,  90, 45, 97, 15),   # row 6
    @(601,  9,  60, 67, 60, 42),   # row 7
    @(1201, 2,  10, 10, 10, 10),   # row 8
    @(1202, 2,  10, 10, 10, 10),   # row 9
    @(902,  1,  0,  0,  0,  0),    # row 10
    @(301,  6,  45, 30, 60, 45),   # row 11
    @(1203, 3,  15, 15, 15, 15),   # row 12
    @(1001, 18, 30, 75, 60, 72),   # row 13
    @(501,  9,  52, 30, 75, 45),   # row 14
    @(801,  3,  67, 65, 52, 45),   # row 15
    @(1101, 0,  15, 30, 30, 0),    # row 16
    @(3,    0,  3,  3,  3,  3),    # row 17
    @(802,  0,  4,  5,  4,  0),    # row 18
    @(2,    0,  2,  2,  2,  2),    # row 19
    @(1,    0,  2,  2,  2,  2)     # row 20
)

$startRow = 3
for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $startRow + $i
    $rowVals = $data[$i]
    for ($c = 0; $c -lt $rowVals.Length; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $rowVals[$c]
    }
}
